$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkmark = [string][char]0x2713

# Row 4: new date header for column I (copy style from H4)
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").Value = 44328

# Rows 5-15: copy style from H column and set checkmark value, mirroring existing pattern
$ws.Range("H5").Copy() | Out-Null
$ws.Range("I5").PasteSpecial(-4122) | Out-Null
$ws.Range("I5").Value = $checkmark

$ws.Range("H6").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Range("I6").Value = $checkmark

$ws.Range("H7").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null
$ws.Range("I7").Value = $checkmark

$ws.Range("G8").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null
$ws.Range("I8").Value = $checkmark

$ws.Range("H9").Copy() | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null
$ws.Range("I9").Value = $checkmark

$ws.Range("H10").Copy() | Out-Null
$ws.Range("I10").PasteSpecial(-4122) | Out-Null
$ws.Range("I10").Value = $checkmark

$ws.Range("H11").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$ws.Range("I11").Value = $checkmark

$ws.Range("H12").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Value = $checkmark

$ws.Range("H13").Copy() | Out-Null
$ws.Range("I13").PasteSpecial(-4122) | Out-Null
$ws.Range("I13").Value = $checkmark

$ws.Range("G14").Copy() | Out-Null
$ws.Range("I14").PasteSpecial(-4122) | Out-Null
$ws.Range("I14").Value = $checkmark

$ws.Range("H15").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = $checkmark

# New empty cell J14, style copied from D5 (style 4, centered, no fill)
$ws.Range("D5").Copy() | Out-Null
$ws.Range("J14").PasteSpecial(-4122) | Out-Null
$ws.Range("J14").ClearContents()

# Update selection to E8
$ws.Range("E8").Select() | Out-Null
